# Update correlation results with JaTyC fixed on DS 6
$wb = $excel.ActiveWorkbook

$allTools = $wb.Worksheets.Item("all_tools")
$typestateChecker = $wb.Worksheets.Item("typestate_checker")

# --- all_tools sheet: rows 10-12 (dataset_id 6) ---
$allTools.Range("F10").Value = 48
$allTools.Range("G10").Value = 819
$allTools.Range("I10").Value = -0.01391413642584883
$allTools.Range("J10").Value = 0.9090032366864784
$allTools.Range("K10").Value = -0.01026612600739247
$allTools.Range("L10").Value = 0.9435901857829745

$allTools.Range("F11").Value = 48
$allTools.Range("G11").Value = 819
$allTools.Range("I11").Value = -0.02815294687312959
$allTools.Range("J11").Value = 0.8086800956106934
$allTools.Range("K11").Value = -0.0299893827986723
$allTools.Range("L11").Value = 0.8362130272367809

$allTools.Range("F12").Value = 48
$allTools.Range("G12").Value = 819
$allTools.Range("I12").Value = -0.247884667760308
$allTools.Range("J12").Value = 0.03188792894690582
$allTools.Range("K12").Value = -0.3096184720415943
$allTools.Range("L12").Value = 0.02866696147428332

# --- typestate_checker sheet: rows 10-12 (dataset_id 6) ---
$typestateChecker.Range("F10").Value = 40
$typestateChecker.Range("G10").Value = 537
$typestateChecker.Range("I10").Value = -0.01391413642584883
$typestateChecker.Range("J10").Value = 0.9090032366864784
$typestateChecker.Range("K10").Value = -0.01026612600739247
$typestateChecker.Range("L10").Value = 0.9435901857829745

$typestateChecker.Range("F11").Value = 40
$typestateChecker.Range("G11").Value = 537
$typestateChecker.Range("I11").Value = -0.02815294687312959
$typestateChecker.Range("J11").Value = 0.8086800956106934
$typestateChecker.Range("K11").Value = -0.0299893827986723
$typestateChecker.Range("L11").Value = 0.8362130272367809

$typestateChecker.Range("F12").Value = 40
$typestateChecker.Range("G12").Value = 537
$typestateChecker.Range("I12").Value = -0.247884667760308
$typestateChecker.Range("J12").Value = 0.03188792894690582
$typestateChecker.Range("K12").Value = -0.3096184720415943
$typestateChecker.Range("L12").Value = 0.02866696147428332

# Column widths widened slightly on typestate_checker sheet (columns I and K)
# Target stored width is 21.7109375 character-units; Excel's column-width
# model only persists pixel-quantised widths, so request the nearest
# reachable ColumnWidth (rounds to a stored width of ~21.667, the closest
# achievable value to 21.7109375).
$typestateChecker.Range("I1").ColumnWidth = 20.8333333333333
$typestateChecker.Range("K1").ColumnWidth = 20.8333333333333
